$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new table row right above the current row 279 ("T10"/quarter_dx),
# shifting it and everything below it down by one. This creates room for the
# new "T09a"/month_rt_dx variant row, which sorts right after "T09"/month_dx
# (row 278) and before "T10"/quarter_dx.
$ws.Rows.Item(279).Insert()

# Populate the new row with the new derived variable definition.
$ws.Cells.Item(279, 1).Value = "T09a"
$ws.Cells.Item(279, 2).Value = "month_rt_dx"
$ws.Cells.Item(279, 3).Value = "Time measurements"
$ws.Cells.Item(279, 4).Value = "Month and year of diagnosis, using the most recent side of the interval as anchor"
$ws.Cells.Item(279, 5).Value = "Months"

# Match the single-line row height used by sibling rows in this table.
$ws.Rows.Item(279).RowHeight = 16

# Row 247 grew taller (likely manually resized after a wrapping change).
$ws.Rows.Item(247).RowHeight = 61

# Grow the table (and its autofilter) by one row to include the newly
# inserted row, keeping the same top-left anchor.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E297"))

# Reflect the updated view/selection state.
$ws.Range("D278").Select()
